$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 168.1098273333333
$ws.Range("N2").Value = 504.329482
$ws.Range("O2").Value = 0.2984182258032519
$ws.Range("P2").Value = 0.298418225803252
$ws.Range("Q2").Value = 2.375279787001778
$ws.Range("R2").Value = 21.377518083016
$ws.Range("S2").Value = 0.2984182258032519
$ws.Range("T2").Value = 0.298418225803252

# Row 3
$ws.Range("O3").Value = 0.2893586437755394
$ws.Range("P3").Value = 0.2893586437755394
$ws.Range("S3").Value = 0.2893586437755394
$ws.Range("T3").Value = 0.2893586437755394

# Row 4
$ws.Range("M4").Value = 165.99353
$ws.Range("N4").Value = 497.98059
$ws.Range("O4").Value = 0.294661504941043
$ws.Range("P4").Value = 0.294661504941043
$ws.Range("Q4").Value = 2.345377916546667
$ws.Range("R4").Value = 21.10840124892
$ws.Range("S4").Value = 0.294661504941043
$ws.Range("T4").Value = 0.294661504941043

# Row 5
$ws.Range("M5").Value = 66.22673433333334
$ws.Range("N5").Value = 198.680203
$ws.Range("O5").Value = 0.1175616254801657
$ws.Range("P5").Value = 0.1175616254801657
$ws.Range("Q5").Value = 0.9357396049737779
$ws.Range("R5").Value = 8.421656444764
$ws.Range("S5").Value = 0.1175616254801657
$ws.Range("T5").Value = 0.1175616254801657
